$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '46.130.09'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.364.58'
$ws.Range("E3").Value = '  +2.33%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '301.99'
$ws.Range("E5").Value = '  +1.08%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '99.76'
$ws.Range("E6").Value = '  +0.60%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.570'
$ws.Range("E7").Value = '  -0.53%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.513'
$ws.Range("E9").Value = '  -2.70%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.45'
$ws.Range("E10").Value = '  -3.79%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0798'
$ws.Range("E11").Value = '  -0.07%  '
$ws.Range("E12").Value = '  -2.50%  '
$ws.Range("E13").Value = '  -0.42%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.728.00'
$ws.Range("E14").Value = '  +2.44%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.388.39'
$ws.Range("E15").Value = '  +3.39%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.815'
$ws.Range("E16").Value = '  -0.09%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.64'
$ws.Range("E17").Value = '  -2.41%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '46.076.88'
$ws.Range("E18").Value = '  -1.59%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.78'
$ws.Range("E19").Value = '  -2.82%  '
$ws.Range("E20").Value = '  +2.81%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.04'
$ws.Range("E21").Value = '  -1.67%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.58'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '245.72'
$ws.Range("E24").Value = '  -2.74%  '
$ws.Range("E25").Value = '  +0.09%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.93'
$ws.Range("E26").Value = '  -2.37%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '39.96'
$ws.Range("E27").Value = '  -6.23%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.20'
$ws.Range("E28").Value = '  -2.84%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.80'
$ws.Range("E29").Value = '  -1.06%  '
$ws.Range("E30").Value = '  +21.98%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '21.03'
$ws.Range("E31").Value = '  +4.08%  '
$ws.Range("E32").Value = '  +6.44%  '
$ws.Range("E33").Value = '  -4.09%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '146.46'
$ws.Range("E34").Value = '  -0.39%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0777'
$ws.Range("E36").Value = '  -1.36%  '
$ws.Range("E37").Value = '  +5.84%  '
$ws.Range("E38").Value = '  -2.28%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '14.96'
$ws.Range("E39").Value = '  -5.31%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.95'
$ws.Range("E40").Value = '  -1.73%  '
$ws.Range("E41").Value = '  -1.88%  '
$ws.Range("E42").Value = '  -5.50%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.911.42'
$ws.Range("E43").Value = '  +3.92%  '
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.11'
$ws.Range("E45").Value = '  +1.83%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.81'
$ws.Range("E46").Value = '  -9.30%  '
$ws.Range("E47").Value = '  -5.94%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.36'
$ws.Range("E48").Value = '  +5.20%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '98.09'
$ws.Range("E49").Value = '  +0.97%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.598.88'
$ws.Range("E50").Value = '  +2.28%  '
$ws.Range("B51").Value = 'ordi'
$ws.Range("C51").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '69.10'
$ws.Range("E51").Value = '  -9.04%  '
